$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update EC2 cost (keep as plain text, not currency number)
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "$20.00"

# Insert new row 7 for ALB, shifting the Total row down to row 9
$ws.Rows("7").Insert()

$ws.Range("A7").Value = "ALB"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "$18.00"

# Update the total (now on row 9)
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "$51.02"
